$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.00002074986032285508
$ws.Range("C2").Value = 0.00007097389502863649
$ws.Range("D2").Value = 3.90043068020848915367
$ws.Range("E2").Value = 8.66023248594897410158
$ws.Range("F2").Value = 0.00000000000000000000
$ws.Range("G2").Value = 12.56075488991280941775

# Row 3
$ws.Range("B3").Value = 0.01514828764759745990
$ws.Range("C3").Value = 114.82701600965050658942
$ws.Range("D3").Value = 49627605961.23487091064453125000
$ws.Range("E3").Value = 9353990175.93243789672851562500
$ws.Range("F3").Value = 0.00000000000000000000
$ws.Range("G3").Value = 58981596252.00948333740234375000
